$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Backfill Arizona (col E / col 5) for 04-10 May 2020, rows 95-101
# ---------------------------------------------------------------------------
$arizonaBackfill = @{
  95 = 0.10633903133903
  96 = 0.1143475572047
  97 = 0.12175324675325
  98 = 0.14635225885226
  99 = 0.13388888888889
  100 = 0.20477855477855
  101 = 0.1404173312068
}
foreach ($r in $arizonaBackfill.Keys) {
  $ws.Cells.Item($r, 5).Value = $arizonaBackfill[$r]
}

# ---------------------------------------------------------------------------
# 2) New date labels for rows 194-195 (column A). Rows 189-193 already have
#    their date label (06-10 Aug 2020); only the data columns are missing.
# ---------------------------------------------------------------------------
$newDateLabels = @{
  194 = "11 08 2020"
  195 = "12 08 2020"
}
foreach ($r in $newDateLabels.Keys) {
  $ws.Cells.Item($r, 1).Value = $newDateLabels[$r]
}

# ---------------------------------------------------------------------------
# 3) Fill in full data for rows 189-193 (06-10 Aug 2020, previously date-only)
#    and new rows 194-195 (11-12 Aug 2020). Column E (Arizona) is left blank
#    for these rows, matching the source data.
# ---------------------------------------------------------------------------
$row189 = @{
  2 = 0.08107899630349499
  3 = 0.11591429841973
  4 = 0.10499700295667
  6 = 0.06887550555208601
  7 = 0.075348475069654
  8 = 0.085988401489462
  9 = 0.074651434986786
  10 = 0.055851976182784
  11 = 0.065353484644079
  12 = 0.07415608878815801
  13 = 0.09335166584479
  14 = 0.082456968819814
  15 = 0.090581075856161
  16 = 0.09146998571689501
  17 = 0.095975975266429
  18 = 0.074689417529529
  19 = 0.087860451557136
  20 = 0.09234757097002699
  21 = 0.083343419108491
  22 = 0.09689934246394701
  23 = 0.072183516910098
  24 = 0.060639344669037
  25 = 0.07392234019983999
  26 = 0.079433118796019
  27 = 0.079049665267831
  28 = 0.093974441825938
  29 = 0.09253209669876999
  30 = 0.11926218367223
  31 = 0.10499768234359
  32 = 0.089449236967587
  33 = 0.09868705232877401
  34 = 0.10381830039617
  35 = 0.081550576068514
  36 = 0.07312551420507001
  37 = 0.088758976723355
  38 = 0.083284006574999
  39 = 0.07078860737254999
  40 = 0.08191574338228
  41 = 0.10683028149455
  42 = 0.082329418692255
  43 = 0.07383805025378901
  44 = 0.088430532391334
  45 = 0.081220351725054
  46 = 0.09762673995646801
  47 = 0.10583445398543
  48 = 0.098561659501007
  49 = 0.091827724410374
  50 = 0.09887131283583001
  51 = 0.077901425107677
  52 = 0.090305908431562
  53 = 0.077783849560894
  54 = 0.079495453059255
  55 = 0.086084781057856
  56 = 0.084057544536881
  57 = 0.11410987679184
}
foreach ($c in $row189.Keys) {
  $ws.Cells.Item(189, $c).Value = $row189[$c]
}

$row190 = @{
  2 = 0.08528850344887499
  3 = 0.099239772539118
  4 = 0.09555880300099701
  6 = 0.07199904048404999
  7 = 0.07419323521774
  8 = 0.090020005537853
  9 = 0.074621911022761
  10 = 0.058001083898971
  11 = 0.068055908782873
  12 = 0.074878182232584
  13 = 0.094411355066793
  14 = 0.078888841425981
  15 = 0.09552851384579
  16 = 0.089401531729426
  17 = 0.09739233397621699
  18 = 0.08134413886118699
  19 = 0.09078419768699
  20 = 0.08700737178451499
  21 = 0.07860391165113299
  22 = 0.091898104617437
  23 = 0.059013424242531
  24 = 0.057035348336114
  25 = 0.066586638717731
  26 = 0.06586909132287801
  27 = 0.059890897741505
  28 = 0.075349127401017
  29 = 0.098139868921782
  30 = 0.1046121447263
  31 = 0.077058441366436
  32 = 0.080020838423872
  33 = 0.085061970378634
  34 = 0.085947614900447
  35 = 0.061323427514395
  36 = 0.053874737585808
  37 = 0.071947928756247
  38 = 0.064442951976314
  39 = 0.055882802136221
  40 = 0.07466829006351899
  41 = 0.095376007411062
  42 = 0.075515794264438
  43 = 0.064899607102971
  44 = 0.077400798097061
  45 = 0.068104388247183
  46 = 0.086847959632756
  47 = 0.09054887556232
  48 = 0.09183917965291299
  49 = 0.083252494209392
  50 = 0.09176723757275999
  51 = 0.069553851258102
  52 = 0.07387658758127599
  53 = 0.065359134799677
  54 = 0.06999146400527
  55 = 0.073845717217437
  56 = 0.07362619903968699
  57 = 0.10076195573437
}
foreach ($c in $row190.Keys) {
  $ws.Cells.Item(190, $c).Value = $row190[$c]
}

$row191 = @{
  2 = 0.07057384242779401
  3 = 0.073446479222755
  4 = 0.066607771746537
  6 = 0.069521365704694
  7 = 0.06480189556687301
  8 = 0.074748897821748
  9 = 0.063428608119501
  10 = 0.057709921654755
  11 = 0.061210528117067
  12 = 0.056709629974563
  13 = 0.055275243638132
  14 = 0.072297209963868
  15 = 0.049415897399954
  16 = 0.062477127153026
  17 = 0.065388352390612
  18 = 0.053352267496902
  19 = 0.058684003656201
  20 = 0.057243974815571
  21 = 0.057666876250653
  22 = 0.06270173081998499
  23 = 0.058582450985422
  24 = 0.050740117656968
  25 = 0.060752848863968
  26 = 0.060333206625039
  27 = 0.056133631807797
  28 = 0.061057372533341
  29 = 0.08694098571525199
  30 = 0.067673132786791
  31 = 0.07229645317589301
  32 = 0.056110170655487
  33 = 0.058062509761545
  34 = 0.06339915759152499
  35 = 0.054823705133703
  36 = 0.051290988224175
  37 = 0.061379017210957
  38 = 0.060827223949659
  39 = 0.050046952281118
  40 = 0.056554373982447
  41 = 0.059225601095127
  42 = 0.059214621469214
  43 = 0.054141005482173
  44 = 0.07105292312801301
  45 = 0.055396485272249
  46 = 0.059644304399622
  47 = 0.058097117591717
  48 = 0.057279370171642
  49 = 0.053523146959504
  50 = 0.068075847350947
  51 = 0.050478759956058
  52 = 0.067809418421363
  53 = 0.053701786930462
  54 = 0.054211754347248
  55 = 0.054997394381643
  56 = 0.055645618642697
  57 = 0.069840875582907
}
foreach ($c in $row191.Keys) {
  $ws.Cells.Item(191, $c).Value = $row191[$c]
}

$row192 = @{
  2 = 0.055724591095004
  3 = 0.062245607100191
  4 = 0.055128404562727
  6 = 0.048723777869942
  7 = 0.053785348655417
  8 = 0.061922221934842
  9 = 0.053424594555461
  10 = 0.046026417304982
  11 = 0.050955585801959
  12 = 0.044255913071989
  13 = 0.046043568379635
  14 = 0.076504129280654
  15 = 0.047344160459536
  16 = 0.057282588311118
  17 = 0.06082139324549
  18 = 0.049898149091969
  19 = 0.054829079614674
  20 = 0.054386683180965
  21 = 0.056183645639025
  22 = 0.060639948999142
  23 = 0.063424452762985
  24 = 0.051678639432347
  25 = 0.06438906461914801
  26 = 0.051331610304879
  27 = 0.043952915885469
  28 = 0.050241548338897
  29 = 0.084453251080645
  30 = 0.053877120191426
  31 = 0.061769824874398
  32 = 0.056521630833433
  33 = 0.054217446650585
  34 = 0.049706893222147
  35 = 0.052259935036747
  36 = 0.051255610098195
  37 = 0.0579596039442
  38 = 0.049960538889244
  39 = 0.056685616514887
  40 = 0.059062250874713
  41 = 0.058068266005323
  42 = 0.060485535106094
  43 = 0.057554515467818
  44 = 0.073229481457615
  45 = 0.064819014491977
  46 = 0.054934358410342
  47 = 0.055995011787941
  48 = 0.053257650511558
  49 = 0.04733471282492
  50 = 0.056897057248801
  51 = 0.049964221919162
  52 = 0.06663026903881999
  53 = 0.051638429643581
  54 = 0.051996817157697
  55 = 0.047389795057878
  56 = 0.050646452747949
  57 = 0.05829878569234
}
foreach ($c in $row192.Keys) {
  $ws.Cells.Item(192, $c).Value = $row192[$c]
}

$row193 = @{
  2 = 0.07074734108401901
  3 = 0.10315529353443
  4 = 0.10262219301479
  6 = 0.06274683670938599
  7 = 0.084219328083076
  8 = 0.091326722449651
  9 = 0.078459665731685
  10 = 0.06547679000399501
  11 = 0.07305044054663801
  12 = 0.07984434464922199
  13 = 0.094563819272566
  14 = 0.078646784046428
  15 = 0.084398076447763
  16 = 0.09504785138560599
  17 = 0.090537375417741
  18 = 0.07002882203238101
  19 = 0.087057787026127
  20 = 0.094067245595239
  21 = 0.08350342154696801
  22 = 0.099304724515951
  23 = 0.070627610835358
  24 = 0.06596245587286199
  25 = 0.07474534285767701
  26 = 0.07293180477215901
  27 = 0.068898379725449
  28 = 0.085086315796691
  29 = 0.10306590006039
  30 = 0.11294688344378
  31 = 0.094349574922233
  32 = 0.087015190326841
  33 = 0.096427277290693
  34 = 0.09678079957671
  35 = 0.070814000670355
  36 = 0.06674282599456299
  37 = 0.077878397757615
  38 = 0.070077722609886
  39 = 0.063474304534012
  40 = 0.07689445281337901
  41 = 0.10103120492813
  42 = 0.075969688192657
  43 = 0.069628938490847
  44 = 0.078935703555418
  45 = 0.068328210581241
  46 = 0.094112997727203
  47 = 0.097349108549691
  48 = 0.096230210099554
  49 = 0.088483366385783
  50 = 0.099622797255232
  51 = 0.07359903095643
  52 = 0.074978150680849
  53 = 0.07071719307813699
  54 = 0.073778250198206
  55 = 0.072805022510013
  56 = 0.076800838710472
  57 = 0.10109429318109
}
foreach ($c in $row193.Keys) {
  $ws.Cells.Item(193, $c).Value = $row193[$c]
}

$row194 = @{
  2 = 0.065968637136207
  3 = 0.10257734045455
  4 = 0.10235963858813
  6 = 0.061206977897677
  7 = 0.081811576839861
  8 = 0.095155553988853
  9 = 0.077620889470032
  10 = 0.059043062315908
  11 = 0.07102227309597001
  12 = 0.076933030793499
  13 = 0.093584849544694
  14 = 0.081756363774152
  15 = 0.085984373586284
  16 = 0.091450434081451
  17 = 0.093337293582238
  18 = 0.070891538120467
  19 = 0.08830703472698199
  20 = 0.094802375713805
  21 = 0.08519611510583699
  22 = 0.1002712782702
  23 = 0.071861092083643
  24 = 0.066931600822375
  25 = 0.07989529924439
  26 = 0.07510534068833
  27 = 0.073718317040503
  28 = 0.087674139073541
  29 = 0.099097141242225
  30 = 0.11699921616367
  31 = 0.09827799119463899
  32 = 0.08727088528664199
  33 = 0.096898746010312
  34 = 0.10340658710587
  35 = 0.074220670500807
  36 = 0.068069582567262
  37 = 0.084070517713452
  38 = 0.073917217664025
  39 = 0.06342278793688701
  40 = 0.07329628696910701
  41 = 0.10191284855256
  42 = 0.076070557315013
  43 = 0.07050997061808301
  44 = 0.07813495525705
  45 = 0.078672284958588
  46 = 0.098167602446165
  47 = 0.10925653166175
  48 = 0.10128605728976
  49 = 0.095854937341221
  50 = 0.10622272807157
  51 = 0.08051611719042299
  52 = 0.080261235171151
  53 = 0.080458378022082
  54 = 0.079816771442508
  55 = 0.07950680027832301
  56 = 0.081168943946008
  57 = 0.10759952911945
}
foreach ($c in $row194.Keys) {
  $ws.Cells.Item(194, $c).Value = $row194[$c]
}

$row195 = @{
  2 = 0.08336663736240001
  3 = 0.10904951748915
  4 = 0.10798981511461
  6 = 0.074480206263128
  7 = 0.074248461816115
  8 = 0.09276439303921501
  9 = 0.07525249414176199
  10 = 0.055501642198015
  11 = 0.07444955080577501
  12 = 0.080430528236246
  13 = 0.094429121062074
  14 = 0.08681421545378901
  15 = 0.08911082669386999
  16 = 0.09657899414608399
  17 = 0.094337234299675
  18 = 0.072518876250192
  19 = 0.093310094860151
  20 = 0.097982719731585
  21 = 0.08970686277222401
  22 = 0.10248101172712
  23 = 0.071216377883305
  24 = 0.067345684469673
  25 = 0.07833128733380899
  26 = 0.07622388133371701
  27 = 0.072037290488479
  28 = 0.08993251118025
  29 = 0.09661787789544
  30 = 0.11508140648987
  31 = 0.09587145529963199
  32 = 0.093995157714217
  33 = 0.10178868736528
  34 = 0.10707359913244
  35 = 0.07727754499461
  36 = 0.070310417782148
  37 = 0.08765623050305101
  38 = 0.07597980074051799
  39 = 0.06866872474100499
  40 = 0.079829338424067
  41 = 0.10680307398489
  42 = 0.077969732020312
  43 = 0.071024867288851
  44 = 0.08305439129499
  45 = 0.071006212687402
  46 = 0.09619130836717001
  47 = 0.10064866379516
  48 = 0.09840847930989501
  49 = 0.08924887037575401
  50 = 0.09877492273593
  51 = 0.074619086075933
  52 = 0.083930453836414
  53 = 0.07222480953435099
  54 = 0.074685261647728
  55 = 0.079851506754476
  56 = 0.076004341395727
  57 = 0.10780648786003
}
foreach ($c in $row195.Keys) {
  $ws.Cells.Item(195, $c).Value = $row195[$c]
}

# ---------------------------------------------------------------------------
# 4) New placeholder rows 196-200 (13-17 Aug 2020): date label only
# ---------------------------------------------------------------------------
$newDateRows = @{
  196 = "13 08 2020"
  197 = "14 08 2020"
  198 = "15 08 2020"
  199 = "16 08 2020"
  200 = "17 08 2020"
}
foreach ($r in $newDateRows.Keys) {
  $ws.Cells.Item($r, 1).Value = $newDateRows[$r]
}

Write-Host "Edit complete."